# Jim's fixes to Module 10 slides
#
# 1) Bump the cached "11/9/2015" auto-date field text (footer/date
#    placeholders) to "11/15/2015" everywhere it is cached: the slide
#    master, every slide layout, and the notes master.
# 2) Rewrite the "Ball Factory has no other behavior" comment + the
#    after-tick/after-button-down/after-button-up/after-drag bodies on
#    slide 10, and reposition + re-flow the callout rectangle beside it.
# 3) Merge the "Study 10-4-ball-factory.rkt..." runs on slide 15 (no
#    text change, just a tidy-up the author's editor performed).

$p = $ppt.ActivePresentation
$d = $p.Designs.Item(1)

# ---------------------------------------------------------------------
# 1) Date placeholders: 11/9/2015 -> 11/15/2015
# ---------------------------------------------------------------------

# Slide master
for ($j = 1; $j -le $d.SlideMaster.Shapes.Count; $j++) {
    $sh = $d.SlideMaster.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = "11/15/2015"
    }
}

# Every slide layout
for ($i = 1; $i -le $d.SlideMaster.CustomLayouts.Count; $i++) {
    $lay = $d.SlideMaster.CustomLayouts.Item($i)
    for ($j = 1; $j -le $lay.Shapes.Count; $j++) {
        $sh = $lay.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "11/15/2015"
        }
    }
}

# Notes master
for ($j = 1; $j -le $p.NotesMaster.Shapes.Count; $j++) {
    $sh = $p.NotesMaster.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = "11/15/2015"
    }
}

# ---------------------------------------------------------------------
# 2) Slide 10 ("Now we can build a ball factory")
# ---------------------------------------------------------------------

$s10 = $p.Slides.Item(10)
$code = $s10.Shapes.Item(2)   # "Content Placeholder 2"
$tr = $code.TextFrame.TextRange

# Paragraph 20: "    ;; the Ball Factory has no other behavior"
#   -> "    ;; the Ball Factory has no other behavior. Return nonsense values for Void,"
#   plus a brand-new paragraph right after it:
#   "    ;; to aid in debugging."
$para20 = $tr.Paragraphs(20, 1)
$para20.Text = "    ;; the Ball Factory has no other behavior. Return nonsense values for Void,"
$para20.InsertAfter("`r    ;; to aid in debugging.")

# Paragraphs 23-26 (after the inserted paragraph + the existing blank
# line): swap the "this" result for the new dummy return values.
$tr.Paragraphs(23, 1).Text = "    (define/public (after-tick) 15)"
$tr.Paragraphs(24, 1).Text = "    (define/public (after-button-down mx my) 16)"
$tr.Paragraphs(25, 1).Text = "    (define/public (after-button-up mx my) 17)"
$tr.Paragraphs(26, 1).Text = "    (define/public (after-drag mx my) 18)"

# The explanatory callout rectangle: move it up/right, next to the code,
# and tidy its run-splitting (text content itself is unchanged).
$rect = $s10.Shapes.Item(4)   # "Rectangle 4"
$rect.Left = 444
$rect.Top = 180

$rtr = $rect.TextFrame.TextRange
$rpara = $rtr.Paragraphs(1, 1)

$lead = 'The factory receives key events from the world.  On each "b", it creates a new ball, and then passes it to the world as an argument to '
$tail = '. '

$leadLen = $lead.Length
$leadRange = $rpara.Characters(1, $leadLen)
$leadRange.Text = "X"
$rpara.Characters(1, 1).Text = $lead

$paraLen = $rpara.Length
$tailRange = $rpara.Characters($paraLen - 1, 2)
$tailRange.Text = "Y"
$rpara.Characters($rpara.Length, 1).Text = $tail

# ---------------------------------------------------------------------
# 3) Slide 15 ("Next Steps"): merge the "Study ..." runs
# ---------------------------------------------------------------------

$s15 = $p.Slides.Item(15)
$content15 = $s15.Shapes.Item(2)   # "Content Placeholder 2"
$para = $content15.TextFrame.TextRange.Paragraphs(1, 1)
$para.Text = "X"
$para.Text = "Study 10-4-ball-factory.rkt in the Examples folder."
